# Updates cryptos list values (price & volume) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.208.55'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.645.36'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('E6').Value = '  +2.12%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +1.17%  '
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.93'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0847'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.876.70'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.651.52'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.14'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('E15').Value = '  +2.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.195.81'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '219.29'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.00%  '
$ws.Range('B22').Value = 'Toncoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.58'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.40%  '
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.20'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '148.03'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.57'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.35%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range('E30').Value = '  -0.75%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('E33').Value = '  +1.38%  '
$ws.Range('E34').Value = '  +1.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.261.96'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.46'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.48%  '
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.544'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.852'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.07%  '
$ws.Range('E40').Value = '  +0.05%  '
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.786.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.80'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.61'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0108'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0513'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.66'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0972'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.10%  '
